# Apply the "double underscore" field separator change for flattened
# fields in the SecondaryMaterialContent reporting template, plus the
# related text/width tweaks described in the commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "semantic_aspect_model_schema": header row 1, columns D..AD
# contain the flattened field names that need "_" -> "__".
# Also widen column D from 49.2 to 50 (to match the other data columns).
# ---------------------------------------------------------------------
$wsSchema = $wb.Worksheets.Item("semantic_aspect_model_schema")

$lastCol = 30  # column AD
for ($col = 4; $col -le $lastCol; $col++) {
    $cell = $wsSchema.Cells.Item(1, $col)
    $current = $cell.Value2
    if ($current -ne $null -and $current.ToString().StartsWith("secondaryMaterialContent[0]_")) {
        $cell.Value = $current.ToString().Replace("_", "__")
    }
}

# Column D width: OOXML width 49.2 -> 50 (ColumnWidth 49.17 renders as 50)
$wsSchema.Columns.Item(4).ColumnWidth = 49.17

# ---------------------------------------------------------------------
# Sheet "description": rows 8..34 in column A hold the same flattened
# field names and need the same "_" -> "__" treatment.
# ---------------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("description")

for ($row = 8; $row -le 34; $row++) {
    $cell = $wsDesc.Cells.Item($row, 1)
    $current = $cell.Value2
    if ($current -ne $null -and $current.ToString().StartsWith("secondaryMaterialContent[0]_")) {
        $cell.Value = $current.ToString().Replace("_", "__")
    }
}

# Legend renumbering: "2. Columns highlighted..." -> "1. Columns highlighted..."
$wsDesc.Range("A3").Value = "1. Columns highlighted in olive green are digital twin fields."

# "Digital Twin Field: x" -> "Digital Twin Field Name: x"
$wsDesc.Range("B5").Value = "Digital Twin Field Name: id"
$wsDesc.Range("B6").Value = "Digital Twin Field Name: manufacturerPartId"
$wsDesc.Range("B7").Value = "Digital Twin Field Name: partInstanceId"
